$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.340.01"
$ws.Range("E2").Value = "  +1.16%  "

$ws.Range("D3").Value = "2.248.50"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "306.69"
$ws.Range("E5").Value = "  -2.21%  "

$ws.Range("D6").Value = "96.43"
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").Value = "35.08"
$ws.Range("E10").Value = "  -1.49%  "

$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  -0.51%  "

$ws.Range("D12").Value = "7.29"
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.342.34"
$ws.Range("E14").Value = "  +4.67%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.591.31"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "0.835"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "13.63"
$ws.Range("E17").Value = "  -2.19%  "

$ws.Range("D18").Value = "44.173.73"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("D19").Value = "0.0₃0969"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").Value = "6.39"
$ws.Range("E20").Value = "  +1.76%  "

$ws.Range("D21").Value = "12.16"
$ws.Range("E21").Value = "  -6.33%  "

$ws.Range("D22").Value = "65.66"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "238.49"
$ws.Range("E23").Value = "  +1.23%  "

$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  -0.91%  "

$ws.Range("D25").Value = "2.01"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "38.68"
$ws.Range("E27").Value = "  +6.23%  "

$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +3.25%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "5.91"
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("D32").Value = "152.39"
$ws.Range("E32").Value = "  -2.81%  "

$ws.Range("D33").Value = "0.0795"
$ws.Range("E33").Value = "  -4.06%  "

$ws.Range("D34").Value = "3.24"
$ws.Range("E34").Value = "  -0.98%  "

$ws.Range("D35").Value = "2.61"
$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("E36").Value = "  +2.77%  "

$ws.Range("E37").Value = "  -2.11%  "

$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -6.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "3.87"
$ws.Range("E40").Value = "  -3.42%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "14.59"
$ws.Range("E41").Value = "  -6.39%  "

$ws.Range("D42").Value = "0.0301"

$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").Value = "1.750.30"
$ws.Range("E44").Value = "  +2.77%  "

$ws.Range("D45").Value = "82.91"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("E46").Value = "  -0.84%  "

$ws.Range("D47").Value = "100.42"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").Value = "4.97"
$ws.Range("E48").Value = "  -2.91%  "

$ws.Range("D49").Value = "8.14"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").Value = "1.58"
$ws.Range("E50").Value = "  -1.29%  "

$ws.Range("D51").Value = "54.91"
$ws.Range("E51").Value = "  -1.96%  "
